$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "51÷7="
$t.Cell(1, 2).Range.Text = "83÷6="
$t.Cell(1, 3).Range.Text = "67÷6="
$t.Cell(1, 4).Range.Text = "80÷5="
$t.Cell(1, 5).Range.Text = "56÷8="
$t.Cell(5, 1).Range.Text = "54÷3="
$t.Cell(5, 2).Range.Text = "67÷5="
$t.Cell(5, 3).Range.Text = "32÷8="
$t.Cell(5, 4).Range.Text = "98÷8="
$t.Cell(5, 5).Range.Text = "59÷5="
$t.Cell(9, 1).Range.Text = "36÷5="
$t.Cell(9, 2).Range.Text = "59÷6="
$t.Cell(9, 3).Range.Text = "43÷5="
$t.Cell(9, 4).Range.Text = "27÷6="
$t.Cell(9, 5).Range.Text = "95÷8="
$t.Cell(13, 1).Range.Text = "53÷8="
$t.Cell(13, 2).Range.Text = "60÷8="
$t.Cell(13, 3).Range.Text = "25÷4="
$t.Cell(13, 4).Range.Text = "13÷5="
$t.Cell(13, 5).Range.Text = "58÷4="
$t.Cell(17, 1).Range.Text = "46÷7="
$t.Cell(17, 2).Range.Text = "87÷8="
$t.Cell(17, 3).Range.Text = "17÷5="
$t.Cell(17, 4).Range.Text = "81÷7="
$t.Cell(17, 5).Range.Text = "71÷7="
